# Add Bird and in Bonus U as referencedTileType
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "W" (Wall/Bird placeholder) tile markers
$ws.Range("AK8").Value = "W"
$ws.Range("AK12").Value = "W"
$ws.Range("AK13").Value = "W"
$ws.Range("AK14").Value = "W"
$ws.Range("AK15").Value = "W"
$ws.Range("AK16").Value = "W"
$ws.Range("AK17").Value = "W"

# New "B" (Bonus) tile markers
$ws.Range("K10").Value = "B"
$ws.Range("AJ10").Value = "B"
$ws.Range("Y11").Value = "B"
$ws.Range("S12").Value = "B"

# New "U" (referenced tile type) markers
$ws.Range("AH17").Value = "U"
$ws.Range("AI17").Value = "U"
$ws.Range("AJ17").Value = "U"

# Update selection to match the authored state
$ws.Range("K10").Select() | Out-Null
